# The workbook contains a header row (row 1) followed by three data rows
# (rows 2-4, columns A:AY). The edit rotates the three data rows "up" by
# one position, with the first data row wrapping around to the end:
#   new row 2 <- old row 3
#   new row 3 <- old row 4
#   new row 4 <- old row 2
# (row 1, the header, is left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"

# Capture the current values of the three data rows before overwriting
# anything, so the rotation can be applied safely.
$row2 = $ws.Range("A2:$($lastCol)2").Value()
$row3 = $ws.Range("A3:$($lastCol)3").Value()
$row4 = $ws.Range("A4:$($lastCol)4").Value()

# Columns I, Y and AA hold values that look like plain numbers / dates
# ("1", "1996-08-10", ...) but are actually stored as literal text in the
# source data. Force those columns to Text formatting before writing the
# values back, so Excel does not "helpfully" reinterpret them as numbers
# or dates during the reassignment below.
$ws.Range("I2:I4").NumberFormat = "@"
$ws.Range("Y2:Y4").NumberFormat = "@"
$ws.Range("AA2:AA4").NumberFormat = "@"

# Write the rotated data back.
$ws.Range("A2:$($lastCol)2").Value = $row3
$ws.Range("A3:$($lastCol)3").Value = $row4
$ws.Range("A4:$($lastCol)4").Value = $row2
